$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("Sheet1")

# Final selection on Sheet1 becomes B4 (and it is no longer the active tab).
$sheet1.Range("B4").Select() | Out-Null

# Add the new "Api" worksheet after Sheet1. We build it by copying Sheet1 so
# that the column widths shared with Sheet1 (columns B/C/D) are carried over
# with full numeric precision instead of being recomputed via AutoFit/ColumnWidth.
$sheet1.Copy([System.Reflection.Missing]::Value, $sheet1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Api"
$ws2.Cells.ClearContents()

# Column A on the new sheet needs its own (narrower) width.
$ws2.Columns("A").ColumnWidth = 11.44140625

# Mark the cells that will hold numeric-looking text (IDs such as "6044") as
# Text so they are stored as strings instead of being coerced into numbers.
$ws2.Range("B1:C1").NumberFormat = "@"
$ws2.Range("A2:D2").NumberFormat = "@"
$ws2.Range("B3:D3").NumberFormat = "@"

# Values are written in the same order the original workbook's shared-string
# table was built in, so new entries land at the same indices: API Data01,
# ScenarioID, CompanyID, 6044, 6045, 1739, 0.
$ws2.Range("A2").Value2 = "API Data01"
$ws2.Range("B1").Value2 = "ScenarioID"
$ws2.Range("C1").Value2 = "CompanyID"
$ws2.Range("B2").Value2 = "6044"
$ws2.Range("B3").Value2 = "6045"
$ws2.Range("C2").Value2 = "1739"
$ws2.Range("D2").Value2 = "0"
$ws2.Range("C3").Value2 = "1739"
$ws2.Range("D3").Value2 = "0"

$ws2.PageSetup.Orientation = 1

# Final selection / active sheet ends up on the new "Api" sheet.
$ws2.Range("E3").Select() | Out-Null
$ws2.Activate() | Out-Null
